$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $orig = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $orig
}

Set-TextValue 'D2' '28.040.30'
Set-TextValue 'E2' '  -3.01%  '
Set-TextValue 'D3' '1.901.80'
Set-TextValue 'E3' '  -3.50%  '
Set-TextValue 'E4' '  -1.08%  '
Set-TextValue 'D5' '326.73'
Set-TextValue 'E5' '  +0.11%  '
Set-TextValue 'D7' '0.4623'
Set-TextValue 'E7' '  -4.12%  '
Set-TextValue 'D8' '0.3964'
Set-TextValue 'E8' '  -2.85%  '
Set-TextValue 'D9' '52.04'
Set-TextValue 'E9' '  -3.45%  '
Set-TextValue 'D10' '0.08314'
Set-TextValue 'E10' '  -3.84%  '
Set-TextValue 'D11' '1.041'
Set-TextValue 'E11' '  -2.39%  '
Set-TextValue 'D12' '21.85'
Set-TextValue 'E12' '  -3.22%  '
Set-TextValue 'D13' '1.932.18'
Set-TextValue 'E13' '  -1.29%  '
Set-TextValue 'D14' '7.364'
Set-TextValue 'E14' '  -5.09%  '
Set-TextValue 'D15' '6.020'
Set-TextValue 'E15' '  -4.25%  '
Set-TextValue 'D16' '1.005'
Set-TextValue 'E16' '  -0.98%  '
Set-TextValue 'E17' '  -2.15%  '
Set-TextValue 'D18' '0.00001061'
Set-TextValue 'E18' '  -1.29%  '
Set-TextValue 'D19' '0.06596'
Set-TextValue 'E19' '  -0.56%  '
Set-TextValue 'D20' '17.78'
Set-TextValue 'E20' '  -5.79%  '
Set-TextValue 'D21' '1.003'
Set-TextValue 'E21' '  -0.88%  '
Set-TextValue 'D22' '5.696'
Set-TextValue 'E22' '  -2.10%  '
Set-TextValue 'D23' '28.054.05'
Set-TextValue 'E23' '  -3.01%  '
Set-TextValue 'D24' '11.13'
Set-TextValue 'E24' '  -4.27%  '
Set-TextValue 'D25' '2.310'
Set-TextValue 'E25' '  +0.88%  '
Set-TextValue 'D26' '2.154.41'
Set-TextValue 'E26' '  -1.51%  '
Set-TextValue 'D27' '153.38'
Set-TextValue 'E27' '  -0.58%  '
Set-TextValue 'D28' '19.96'
Set-TextValue 'E28' '  -2.22%  '
Set-TextValue 'D29' '2.120'
Set-TextValue 'E29' '  -2.80%  '
Set-TextValue 'D30' '5.706'
Set-TextValue 'E30' '  -5.94%  '
Set-TextValue 'D31' '123.94'
Set-TextValue 'E31' '  -0.98%  '
Set-TextValue 'B32' 'Stellar'
Set-TextValue 'C32' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D32' '0.09581'
Set-TextValue 'E32' '  -1.01%  '
Set-TextValue 'B33' 'ImmutableX'
Set-TextValue 'C33' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D33' '0.9672'
Set-TextValue 'E33' '  -4.70%  '
Set-TextValue 'E34' '  -0.97%  '
Set-TextValue 'D35' '3.622'
Set-TextValue 'E35' '  -2.18%  '
Set-TextValue 'D36' '5.498'
Set-TextValue 'E36' '  -4.17%  '
Set-TextValue 'D37' '1.259'
Set-TextValue 'E37' '  -1.98%  '
Set-TextValue 'B38' 'FraxShare'
Set-TextValue 'C38' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D38' '8.702'
Set-TextValue 'E38' '  -1.80%  '
Set-TextValue 'B39' 'VeChain'
Set-TextValue 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D39' '0.02278'
Set-TextValue 'E39' '  -3.67%  '
Set-TextValue 'D40' '0.06123'
Set-TextValue 'D41' '0.6120'
Set-TextValue 'E41' '  -2.78%  '
Set-TextValue 'D42' '1.002'
Set-TextValue 'E42' '  -0.90%  '
Set-TextValue 'D43' '10.83'
Set-TextValue 'E43' '  -3.74%  '
Set-TextValue 'D44' '0.1900'
Set-TextValue 'E44' '  -1.49%  '
Set-TextValue 'D45' '1.300'
Set-TextValue 'E45' '  -2.57%  '
Set-TextValue 'D46' '0.5839'
Set-TextValue 'E46' '  -2.88%  '
Set-TextValue 'E47' '  -2.20%  '
Set-TextValue 'D48' '2.000'
Set-TextValue 'E48' '  -4.89%  '
Set-TextValue 'D49' '3.436'
Set-TextValue 'E49' '  -0.36%  '
Set-TextValue 'D50' '0.06894'
Set-TextValue 'E50' '  +0.35%  '
Set-TextValue 'D51' '110.42'
Set-TextValue 'E51' '  -0.90%  '
